$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Z-column formula: the "-M<row>" term should be "-(12*M<row>)".
# Apply to every row (2 through 25) that holds this formula.
for ($r = 2; $r -le 25; $r++) {
    $formula = "=((0.6+(0.02*(T$r-2005)))*(2*12*M$r))-(12*M$r)+12*(P$r+N$r+0.96*O$r)"
    $ws.Range("Z$r").Formula = $formula
}

# Update the view's selection to match the corrected workbook state.
$ws.Range("P9").Select() | Out-Null
